# Auto-generated script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 630.05
$ws.Range("J38").Value = 3000
$ws.Range("L38").Value = 9000
$ws.Range("N38").Value = -9744
$ws.Range("H58").Value = 92747.27
$ws.Range("I58").Value = 766.6667
$ws.Range("J58").Value = 127240
$ws.Range("K58").Value = 2300.0001
$ws.Range("L58").Value = 381720
$ws.Range("M58").Value = -2150.0001
$ws.Range("N58").Value = -382020
$ws.Range("H87").Value = 22541.586
$ws.Range("J87").Value = 22541.586
$ws.Range("L87").Value = 22541.586
$ws.Range("N87").Value = -25037.586
$ws.Range("H90").Value = 22541.586
$ws.Range("J90").Value = 22541.586
$ws.Range("L90").Value = 67624.758
$ws.Range("N90").Value = -80104.758
$ws.Range("H97").Value = 6490
$ws.Range("J97").Value = 6490
$ws.Range("L97").Value = 19470
$ws.Range("N97").Value = -20462
$ws.Range("H99").Value = 1296.2142
$ws.Range("J99").Value = 1975.125
$ws.Range("L99").Value = 5925.375
$ws.Range("N99").Value = -8921.375
$ws.Range("H112").Value = 9616328
$ws.Range("J112").Value = 10000961
$ws.Range("L112").Value = 30002883
$ws.Range("N112").Value = -30005099
$ws.Range("H129").Value = 4808765
$ws.Range("I129").Value = 50001360
$ws.Range("J129").Value = 1042.1489
$ws.Range("K129").Value = 150004080
$ws.Range("L129").Value = 3126.4467
$ws.Range("M129").Value = -149999080
$ws.Range("N129").Value = -13126.4467
$ws.Range("H137").Value = 2568182.8
$ws.Range("I137").Value = 3708355.2
$ws.Range("J137").Value = 2794.5833
$ws.Range("K137").Value = 11125065.6
$ws.Range("L137").Value = 8383.749899999999
$ws.Range("M137").Value = -11122515.6
$ws.Range("N137").Value = -13483.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5876.841
$ws.Range("I32").Value = 5044.012
$ws.Range("J32").Value = 19701.8
$ws.Range("K32").Value = 5044.012
$ws.Range("L32").Value = 19701.8
$ws.Range("M32").Value = -4757.012
$ws.Range("N32").Value = -20275.8
$ws.Range("H62").Value = 31830
$ws.Range("I62").Value = 980
$ws.Range("J62").Value = 38000
$ws.Range("K62").Value = 980
$ws.Range("L62").Value = 38000
$ws.Range("M62").Value = -356
$ws.Range("N62").Value = -39248
$ws.Range("H63").Value = 2059.7646
$ws.Range("I63").Value = 1876
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 1876
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1190
$ws.Range("N63").Value = -6372
$ws.Range("H65").Value = 31830
$ws.Range("I65").Value = 980
$ws.Range("J65").Value = 38000
$ws.Range("K65").Value = 2940
$ws.Range("L65").Value = 114000
$ws.Range("M65").Value = 180
$ws.Range("N65").Value = -120240
$ws.Range("H66").Value = 2059.7646
$ws.Range("I66").Value = 1876
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 9380
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -5948
$ws.Range("N66").Value = -31864
$ws.Range("H97").Value = 806.625
$ws.Range("I97").Value = 786.4706
$ws.Range("J97").Value = 855.5714
$ws.Range("K97").Value = 786.4706
$ws.Range("L97").Value = 855.5714
$ws.Range("M97").Value = -290.4706
$ws.Range("N97").Value = -1847.5714
$ws.Range("H132").Value = 4647.294
$ws.Range("I132").Value = 4272
$ws.Range("J132").Value = 5335.3335
$ws.Range("K132").Value = 12816
$ws.Range("L132").Value = 16006.0005
$ws.Range("M132").Value = -10286
$ws.Range("N132").Value = -21066.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 4075.889
$ws.Range("I11").Value = 135.6
$ws.Range("J11").Value = 9001.25
$ws.Range("K11").Value = 135.6
$ws.Range("L11").Value = 9001.25
$ws.Range("M11").Value = 4.400000000000006
$ws.Range("N11").Value = -9281.25
$ws.Range("H134").Value = 7380
$ws.Range("I134").Value = 5633.3335
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 16900.0005
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -14365.0005
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2786
$ws.Range("I22").Value = 965
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 965
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -615
$ws.Range("N22").Value = -4700
$ws.Range("H31").Value = 1889051.9
$ws.Range("I31").Value = 2223815.5
$ws.Range("J31").Value = 6007.375
$ws.Range("K31").Value = 2223815.5
$ws.Range("L31").Value = 6007.375
$ws.Range("M31").Value = -2223520.5
$ws.Range("N31").Value = -6597.375
$ws.Range("H34").Value = 1889051.9
$ws.Range("I34").Value = 2223815.5
$ws.Range("J34").Value = 6007.375
$ws.Range("K34").Value = 2223815.5
$ws.Range("L34").Value = 6007.375
$ws.Range("M34").Value = -2223613.5
$ws.Range("N34").Value = -6411.375
$ws.Range("H52").Value = 26666.666
$ws.Range("J52").Value = 26666.666
$ws.Range("L52").Value = 26666.666
$ws.Range("N52").Value = -27254.666
$ws.Range("H132").Value = 3146
$ws.Range("I132").Value = 2616.6667
$ws.Range("J132").Value = 3940
$ws.Range("K132").Value = 7850.000100000001
$ws.Range("L132").Value = 11820
$ws.Range("M132").Value = -5320.000100000001
$ws.Range("N132").Value = -16880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1160.6
$ws.Range("I131").Value = 2135
$ws.Range("K131").Value = 6405
$ws.Range("M131").Value = -1365
$ws.Range("H132").Value = 3000
$ws.Range("J132").Value = 6500
$ws.Range("L132").Value = 58500
$ws.Range("N132").Value = -63560
$ws.Range("H137").Value = 3327.524
$ws.Range("I137").Value = 4044
$ws.Range("J137").Value = 3230.7026
$ws.Range("K137").Value = 12132
$ws.Range("L137").Value = 9692.1078
$ws.Range("M137").Value = -7032
$ws.Range("N137").Value = -19892.1078

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3328.5
$ws.Range("I43").Value = 2330
$ws.Range("J43").Value = 3927.6
$ws.Range("K43").Value = 2330
$ws.Range("L43").Value = 3927.6
$ws.Range("M43").Value = -2179
$ws.Range("N43").Value = -4229.6
$ws.Range("H132").Value = 4820.56
$ws.Range("I132").Value = 5626.3335
$ws.Range("J132").Value = 4076.7693
$ws.Range("K132").Value = 16879.0005
$ws.Range("L132").Value = 12230.3079
$ws.Range("M132").Value = -14349.0005
$ws.Range("N132").Value = -17290.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5264838
$ws.Range("I7").Value = 9091819
$ws.Range("J7").Value = 2738.125
$ws.Range("K7").Value = 9091819
$ws.Range("L7").Value = 2738.125
$ws.Range("M7").Value = -9091707
$ws.Range("N7").Value = -2962.125
$ws.Range("H126").Value = 5264838
$ws.Range("I126").Value = 9091819
$ws.Range("J126").Value = 2738.125
$ws.Range("K126").Value = 27275457
$ws.Range("L126").Value = 8214.375
$ws.Range("M126").Value = -27272987
$ws.Range("N126").Value = -13154.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 160055.44
$ws.Range("I132").Value = 185935.1
$ws.Range("K132").Value = 557805.3
$ws.Range("M132").Value = -555275.3
